# Bai 8 Danh gia ket qua tim kiem_phan 2.pptx
# Commit: "Them tt chuong sach TA" (Add English-book chapter abbreviation)
#
# The subtitle on slide 1 currently reads:
#   "IIR.Chap8. Evaluation in information retrieval"
# It should read:
#   "IIR.C8. Evaluation in information retrieval"
#
# i.e. "Chap8" -> "C8" inside the run that starts with "IIR.".

$p = $ppt.ActivePresentation

$target = "IIR.Chap8"
$replacement = "IIR.C8"
$done = $false

for ($si = 1; $si -le $p.Slides.Count -and -not $done; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count -and -not $done; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -eq $false) { continue }
        if ($shape.HasTextFrame -ne -1) { continue }

        $tf = $shape.TextFrame
        if ($tf.HasText -eq $false) { continue }

        $tr = $tf.TextRange
        $hit = $tr.Find($target, 0)
        if ($hit -ne $null) {
            $hit.Text = $replacement
            $done = $true
        }
    }
}

if (-not $done) {
    throw "Could not find text '$target' to replace in presentation."
}
